$d = $word.ActiveDocument

# 1. Title: append " Teórica e Prática"
$d.Content.Find.Execute('LOT2040 -  Engenharia Genética', $true, $false, $false, $false, $false, $true, 1, $false, 'LOT2040 -  Engenharia Genética Teórica e Prática', 2) | Out-Null

# 2. Ativação date 2019 -> 2025
$d.Content.Find.Execute('Ativação: 01/01/2019', $true, $false, $false, $false, $false, $true, 1, $false, 'Ativação: 01/01/2025', 2) | Out-Null

# 3. Curso (semestre ideal): EB (2) -> EB (6)
$d.Content.Find.Execute('Curso (semestre ideal): EB (2)', $true, $false, $false, $false, $false, $true, 1, $false, 'Curso (semestre ideal): EB (6)', 2) | Out-Null

# 4. Fill the empty italic English "Objetivos" run with its text
#    (the paragraph right after the Portuguese objectives paragraph, whose only
#     run currently has italic rPr but no text)
$objParaIdx = 0
for ($i = 2; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    $prev = $d.Paragraphs.Item($i - 1)
    if (($pp.Range.Text.Trim() -eq "") -and ($prev.Range.Text -like "Introduzir os princípios básicos*")) {
        $objParaIdx = $i
    }
}
if ($objParaIdx -gt 0) {
    $d.Paragraphs.Item($objParaIdx).Range.Text = 'Introduce fundamental genetic principles while thoroughly covering essential concepts in classical genetics, molecular genetics, and genomics, alongside the techniques and applications of molecular genetics.'
}

# 5. English program summary (italic) full replacement
$d.Content.Find.Execute('Introduction to Genetics; Structure and Replication of DNA; Transcription, Translation and the Genetic Code; Mutation and DNA Repair; Regulation of Gene Expression; Genomics and Bioinformatics; The Techniques and Applications of Molecular Genetics', $true, $false, $false, $false, $false, $true, 1, $false, 'Introduction to genetics; Structure and Molecular Replication of DNA; Transcription, Translation and Genetic Code; DNA Mutation and Repair; Regulation of Gene Expression; Genomics and Bioinformatics; Techniques and Applications of Molecular Genetics.', 2) | Out-Null

# 6. Portuguese numbered "Programa" list spacing/punctuation
$d.Content.Find.Execute('1. Introdução à Genética.2. Nucleotídeos e Estrutura do DNA3. Genes e Cromossomos4. Replicação do DNA5. Transcrição e Processamento do RNA6. Código genético e Tradução7. Clonagem de DNA8. Genômica9. Bioinformática10. Regulação da Expressão Gênica11. Elementos Genéticos Transponíveis12. Mutação, Reparo de DNA e Recombinação13. Técnicas e Aplicações da Genética Molecular', $true, $false, $false, $false, $false, $true, 1, $false, '1. Introdução à Genética. 2. Nucleotídeos e Estrutura do DNA 3. Genes e Cromossomos 4. Replicação do DNA 5. Transcrição e Processamento do RNA 6. Código genético e Tradução 7. Clonagem de DNA 8. Genômica 9. Bioinformática 10. Regulação da Expressão Gênica 11. Elementos Genéticos Transponíveis 12. Mutação, Reparo de DNA e Recombinação 13. Técnicas e Aplicações da Genética Molecular.', 2) | Out-Null

# 7. English numbered "Programa" list spacing/punctuation/rewording
$d.Content.Find.Execute('1. Introduction to Genetics Science.2. Nucleotides and Structure of DNA3. Genes and Chromosomes4. DNA replication5. Transcription and RNA processing6. The Genetic Code and Translation7. DNA cloning8. Genomics9. Bioinformatics10. Regulation of Gene Expression11. Transposable Genetic Elements12. Mutation, DNA Repair, and Recombination13. Techniques and Applications of Molecular Genetics', $true, $false, $false, $false, $false, $true, 1, $false, '1. Introduction to Genetics. 2. Nucleotides and DNA Structure 3. Genes and Chromosomes 4. DNA Replication 5. Transcription and RNA Processing 6. Genetic Code and Translation 7. DNA Cloning 8. Genomics 9. Bioinformatics 10. Regulation of Gene Expression 11. Elements Transposable Genetics 12. Mutation, DNA Repair and Recombination 13. Techniques and Applications of Molecular Genetics', 2) | Out-Null

# 8. Norma de recuperação text (prepend label text + extra space)
$d.Content.Find.Execute('NF = (MF + PR)/2, onde PR é uma prova de recuperação.Prova de Recuperação (PR) para alunos com Média Final (MF) maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final (NF) igual ou maior do que 5,0.', $true, $false, $false, $false, $false, $true, 1, $false, 'Norma de Recuperação: NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de Recuperação (PR) para alunos com Média Final (MF) maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final (NF) igual ou maior do que 5,0.', 2) | Out-Null

# 9. Bibliography full rewrite
$d.Content.Find.Execute('- Fundamentos de Genética. Peter Snustad e Michael Simmons (2013 Sexta Edição). Editora Guanabara Koogan- Principles of Genetics. Dr. Peter Snustad and Michael Simmons (2016 7th Edition). John Wiley and Sons- Genética: Um enfoque conceitual. Benjamin Pierce (2016 – Quinta Edição). Editora Guanabara Koogan.- Genetics: A Conceptual Approach. Benjamin Pierce (2016 – Sixth Edition). W.H. Freeman and Company.- Introdução à genética. Lewontin, Griffiths, Carroll e Wessler. (2013 – Décima Edição). Guanabara Koogan.- An Introduction to Genetic Analysis. Anthony Griffiths, Susan Wessler, Sean Carroll, and John Doebley (2015 Eleventh Edition). W. H. Freeman', $true, $false, $false, $false, $false, $true, 1, $false, '-Fundamentos de Genética. Peter Snustad e Michael Simmons (2017, 7ª Edição). Editora Guanabara Koogan. -Principles of Genetics. Dr. Peter Snustad and Michael Simmons (2019, 7th Edition). John Wiley and Sons. -Genética: Um enfoque conceitual. Benjamin Pierce (2016, 5ª Edição). Editora Guanabara Koogan. -Genetics: A Conceptual Approach. Benjamin Pierce (2019, 7th Edition). W. H. Freeman. -Introdução à genética. Griffiths, Doebley, Peichel e Wassarman (2022 – 12ª Edição). Guanabara Koogan. -An Introduction to Genetic Analysis. Anthony Griffiths, John Doebley, Catherine Peichel, David A. Wassarman (2020 12th Edition). W. H. Freeman.', 2) | Out-Null

# 10. Requisitos: insert new LOT2053 bullet line before LOT2002 line,
#     as its own run (with its own trailing line break) ahead of the
#     existing run, matching the target run structure.
$reqPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$reqRange = $reqPara.Range
$insPoint = $reqRange.Duplicate
$insPoint.Collapse(1)
$lineBreak = [char]11
$newBulletText = 'LOT2053 -  Microbiologia: da Teoria à Prática  (Requisito fraco)' + $lineBreak
$insPoint.InsertBefore($newBulletText)
$newRunRange = $d.Range($reqRange.Start, $reqRange.Start + $newBulletText.Length)
# Toggle direct formatting on/off so the engine commits the inserted text
# as its own run instead of merging it into the following run.
$newRunRange.Bold = 1
$newRunRange.Bold = 0
